# "Generate Report for Archive"
# - Status text moves from "Ready for handoff" to "In Translation" everywhere
#   it is used (Overview!E2:F4, zh-cn!C2:C4, de-de!C2:C4).
# - The Status/zh-cn/de-de columns get narrower to fit the new (shorter) text:
#   Overview columns E:F and the "Status" column (C) on the zh-cn / de-de
#   sheets shrink from ~17.22 chars to ~13.41 chars.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# Target stored column width is ~13.41 characters; ColumnWidth assignments
# are quantized to the nearest pixel by the host, so 12.5 is the input that
# lands closest to that target.
$newColWidth = 12.5

# --- Overview sheet: columns E (zh-cn) and F (de-de), rows 2-4 ---
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($row in 2..4) {
    foreach ($col in @("E", "F")) {
        $cell = $wsOverview.Range($col + $row)
        if ($cell.Text -eq $oldStatus) {
            $cell.Value = $newStatus
        }
    }
}
$wsOverview.Range("E1:F1").ColumnWidth = $newColWidth

# --- zh-cn / de-de sheets: "Status" column C, rows 2-4 ---
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in 2..4) {
        $cell = $ws.Range("C" + $row)
        if ($cell.Text -eq $oldStatus) {
            $cell.Value = $newStatus
        }
    }
    $ws.Range("C1").ColumnWidth = $newColWidth
}
